# Commit: "addind preparer to sheet"
# Replace the placeholder "Retrofitted_3067" value in the libraryPreparer (B)
# and purpose (E) columns with the actual preparer's initials "H.BROWN"
# for all data rows (2-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B11").Value = "H.BROWN"
$ws.Range("E2:E11").Value = "H.BROWN"

# Reflect the selection state left behind in the source workbook
# (active cell E2, within the combined B3:B11 / E2:E11 selection).
$ws.Range("E2:E11,B3:B11").Select()
